# reminder.xlsx — "make changes in the task follow-up task"
#
# 1. Reword the bold lead-in of the "title" note question
#      " **What is the title of the task?**:"  ->  " **Task title**:"
# 2. Reword the bold lead-in of the "when" note question
#      " **When should a reminder for this task appear?**"  ->  " **Date task due?**"
# 3. Remove the "dt" note question ("**Date Task Appears:** ${format}") from the
#    survey sheet entirely — this shifts every row below it up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# --- 1. C27 ("title" question label) -------------------------------------
# Rich text: run 1 is bold (" **What is the title of the task?**:"), run 2 is
# the plain "  ${titt_task} " reference — only run 1's text changes.
$titleCell = $ws.Range("C27")
$oldLead = " **What is the title of the task?**:"
$titleCell.Characters(1, $oldLead.Length).Text = " **Task title**:"

# --- 2. C29 ("when" question label) ---------------------------------------
# Same pattern: run 1 is bold (" **When should a reminder for this task
# appear?**"), run 2 is the plain "  ${format} " reference.
$whenCell = $ws.Range("C29")
$oldLead2 = " **When should a reminder for this task appear?**"
$whenCell.Characters(1, $oldLead2.Length).Text = " **Date task due?**"

# --- 3. Remove the "dt" note row (row 31) ----------------------------------
# type=note, name=dt, label="**Date Task Appears:** ${format} " — deleting the
# whole row shifts the follow-up group (select_one follow / string follow_up /
# end group) up by one row, matching the target layout.
$ws.Rows.Item(31).Delete()
